# Updated labels and fixed jail time credit bug.
# Appends additional case rows (805-816) to Sheet1, mirroring the rows that
# were added to db/Case_Data.xlsx, and extends the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("22CRB00136", "Hemmeter", "DOMESTIC VIOLENCE",                     "2919.25(A)",    "No Data", "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("22CRB00136", "Hemmeter", "ASSAULT - M1",                          "2903.13(A)",    "No Data", "No Contest", "Guilty", "$ 0",  "$ 0",  "None", "None"),
    @("22CRB00136", "Hemmeter", "DOMESTIC VIOLENCE",                     "2919.25(A)",    "No Data", "No Contest", "Guilty", "$ 50", "$ 25", "10",   "None"),
    @("22CRB00136", "Hemmeter", "ASSAULT - M1",                          "2903.13(A)",    "No Data", "No Contest", "Guilty", "$ 0",  "$ 0",  "None", "None"),
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM",      "No Contest", "Guilty", "$ 10", "$ 5",  "None", "None")
)

$startRow = 805
$endRow = $startRow + $newRows.Length - 1

# Force the new range to be plain text so that numeric-looking strings like
# "10" and currency-looking strings like "$ 50" are stored as text instead
# of being auto-converted into numbers by Excel.
$fullRange = $ws.Range("A" + $startRow + ":K" + $endRow)
$fullRange.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $ws.Cells.Item($r, 8).Value = $rowData[7]
    $ws.Cells.Item($r, 9).Value = $rowData[8]
    $ws.Cells.Item($r, 10).Value = $rowData[9]
    $ws.Cells.Item($r, 11).Value = $rowData[10]
}
